# feat: add 2022-Q3 data
#
# The workbook tracks quarterly fund holdings. A new quarter ("2022-Q3")
# is inserted right after the "总计" (summary) sheet and before the
# existing "2022-Q2" sheet, which (together with "2022-Q1") simply shifts
# one position to the right. The summary sheet gets a new row for the
# new quarter, and the previous two rows' labels shift down one quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: relabel the existing two data
#    rows one quarter forward, and append a new row for "2022-Q1" that
#    duplicates the (unchanged) totals, matching the existing pattern.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Clone row 3's formatting down into the new row 4 first (so the new
# row's style matches rows 2/3 exactly), then overwrite values.
$summary.Range("A3:D3").Copy()
$summary.Range("A4:D4").PasteSpecial(-4122)

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("B3").Value = "2022-Q2"

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.01

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" sheet. Copying the existing "2022-Q2"
#    sheet (placing the copy right before it) gives the new sheet the
#    exact same layout/formatting, and pushes "2022-Q2"/"2022-Q1" one
#    tab to the right - matching the target sheet order.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Copy($oldQ2)

$newQ3 = $wb.Worksheets.Item(2)
$newQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 3) Fill in the new quarter's fund figures on "2022-Q3".
# ---------------------------------------------------------------------
$newQ3.Range("C2").Value = "建信新兴市场优选混合（QDII）"

$newQ3.Range("D2").NumberFormat = "@"
$newQ3.Range("D2").Value = "0.14"

$newQ3.Range("E2").NumberFormat = "@"
$newQ3.Range("E2").Value = "81.57"

$newQ3.Range("F2").NumberFormat = "@"
$newQ3.Range("F2").Value = "6.63"

$newQ3.Range("G2").NumberFormat = "@"
$newQ3.Range("G2").Value = "0.0093"

# ---------------------------------------------------------------------
# 4) Keep "2022-Q1" as the selected tab, as it was before this edit.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
